$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Setor de Governo" (column C) for the INFRAESTRUTURA, MOBILIDADE E
# PARCERIAS rows (103-111): previously duplicated the full secretariat name,
# now uses the shorter sector label.
for ($r = 103; $r -le 111; $r++) {
    $ws.Range("C$r").Value = "INFRAESTRUTURA, MOBILIDADE E PARCERIAS"
}

# Widen column C to fit the text.
$ws.Range("C1").ColumnWidth = 46

# Move the active selection to C1.
$ws.Range("C1").Select()
